# Update of "Poland Ekstraklasa" odds base, 02-04-2024 23:59
# - Row 234 (match id 232, Radomiak Radom vs Zaglebie Lubin) gets its
#   result + closing-line columns filled in (it had just kick-off/opening
#   odds before).
# - Five new finished/updated matches are appended as rows 236-240
#   (row 235 is untouched, pre-existing context).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 234 : fill in result (H/I/J) and refresh closing odds columns
# ---------------------------------------------------------------------
$ws.Cells.Item(234, 8).Value2  = 1        # H234 FTHG
$ws.Cells.Item(234, 9).Value2  = 1        # I234 FTAG
$ws.Cells.Item(234, 10).Value2 = "D"      # J234 FTR

$ws.Cells.Item(234, 14).Value2 = 3.1      # N234 oddH
$ws.Cells.Item(234, 16).Value2 = 2.55     # P234 oddA
$ws.Cells.Item(234, 18).Value2 = 2.1      # R234 oddAHH
$ws.Cells.Item(234, 19).Value2 = 1.775    # S234 oddAHA
$ws.Cells.Item(234, 21).Value2 = 1.9      # U234 oddAHOver
$ws.Cells.Item(234, 22).Value2 = 1.95     # V234 oddAHUnder
$ws.Cells.Item(234, 23).Value2 = -1       # W234 PLH
$ws.Cells.Item(234, 24).Value2 = 2        # X234 PLD
$ws.Cells.Item(234, 25).Value2 = -1       # Y234 PLA
$ws.Cells.Item(234, 27).Value2 = -0       # AA234 PL_Ahh

# AB234/AC234 did not exist yet -> copy number formatting from the
# neighbouring AA234 cell before writing the new values.
$ws.Range("AA234").Copy()
$ws.Range("AB234:AC234").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(234, 28).Value2 = 0        # AB234 PL_Aha
$ws.Cells.Item(234, 29).Value2 = -0       # AC234 PL_AhOver

# ---------------------------------------------------------------------
# Append rows 236-240 (row 235 already exists and stays untouched)
# ---------------------------------------------------------------------
function Add-Match {
    param($Row, $Id, $MatchId, $Date, $Home, $Away, $Odds)

    # Clone formatting of only the id (A) / date (E) cells from the row
    # above, so those keep their style (bold/centered id, date number
    # format, ...) without stamping empty-but-styled placeholder cells
    # into the columns that stay blank (H/I/J - no result yet).
    $srcRow = $Row - 1
    $ws.Range("A" + $srcRow).Copy()
    $ws.Range("A" + $Row).PasteSpecial(-4122)
    $ws.Range("E" + $srcRow).Copy()
    $ws.Range("E" + $Row).PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Cells.Item($Row, 1).Value2  = $Id                   # A  id
    $ws.Cells.Item($Row, 2).Value2  = $MatchId               # B
    $ws.Cells.Item($Row, 3).Value2  = "Poland Ekstraklasa"   # C Div
    $ws.Cells.Item($Row, 4).Value2  = "Poland Ekstraklasa"   # D Div Original Name
    $ws.Cells.Item($Row, 5).Value2  = $Date                  # E Date
    $ws.Cells.Item($Row, 6).Value2  = $Home                  # F HomeTeam
    $ws.Cells.Item($Row, 7).Value2  = $Away                  # G AwayTeam
    # H,I,J (FTHG/FTAG/FTR) left blank - match not played yet
    $ws.Cells.Item($Row, 11).Value2 = $Odds[0]   # K oddH_op
    $ws.Cells.Item($Row, 12).Value2 = $Odds[1]   # L oddD_op
    $ws.Cells.Item($Row, 13).Value2 = $Odds[2]   # M oddA_op
    $ws.Cells.Item($Row, 14).Value2 = $Odds[3]   # N oddH
    $ws.Cells.Item($Row, 15).Value2 = $Odds[4]   # O oddD
    $ws.Cells.Item($Row, 16).Value2 = $Odds[5]   # P oddA
    $ws.Cells.Item($Row, 17).Value2 = $Odds[6]   # Q Ah
    $ws.Cells.Item($Row, 18).Value2 = $Odds[7]   # R oddAHH
    $ws.Cells.Item($Row, 19).Value2 = $Odds[8]   # S oddAHA
    $ws.Cells.Item($Row, 20).Value2 = $Odds[9]   # T AhOU
    $ws.Cells.Item($Row, 21).Value2 = $Odds[10]  # U oddAHOver
    $ws.Cells.Item($Row, 22).Value2 = $Odds[11]  # V oddAHUnder
    $ws.Cells.Item($Row, 23).Value2 = $Odds[12]  # W PLH
    $ws.Cells.Item($Row, 24).Value2 = $Odds[13]  # X PLD
    $ws.Cells.Item($Row, 25).Value2 = $Odds[14]  # Y PLA
    $ws.Cells.Item($Row, 26).Value2 = $Odds[15]  # Z PL_Ahh
    $ws.Cells.Item($Row, 27).Value2 = $Odds[16]  # AA PL_Aha
}

Add-Match 236 234 6775585 45387.64583333334 "Radomiak Radom" "Rakow Czestochowa" `
    @(4, 3.5, 1.909, 4.2, 3.6, 1.85, 0.5, 2, 1.85, 2.25, 1.85, 2, 0, 0, 0, 0, 0)

Add-Match 237 235 6775582 45388.41666666666 "Korona Kielce" "Stal Mielec" `
    @(1.8, 3.5, 4.5, 1.75, 3.5, 4.75, -0.75, 2.025, 1.825, 2.25, 1.875, 1.975, 0, 0, 0, 0, 0)

Add-Match 238 236 6774876 45388.52083333334 "Ruch Chorzow" "Puszcza Niepolomice" `
    @(1.85, 3.5, 4.2, 1.8, 3.6, 4.5, -0.75, 2.05, 1.8, 2.25, 1.8, 2.05, 0, 0, 0, 0, 0)

Add-Match 239 237 6775586 45389.3125 "Widzew Lodz" "Piast Gliwice" `
    @(2.5, 3.25, 2.8, 2.625, 3.25, 2.7, 0, 1.925, 1.925, 2.25, 2.05, 1.8, 0, 0, 0, 0, 0)

Add-Match 240 238 6775583 45389.41666666666 "Lech Poznan" "Pogon Szczecin" `
    @(2.5, 3.4, 2.7, 2.45, 3.4, 2.7, 0, 1.825, 2.025, 2.5, 1.825, 2.025, 0, 0, 0, 0, 0)
